# Update the "year" column (A3:A6) on both sheets from 2011-2014 to 2007-2010
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Page1")
$ws1.Range("A3").Value = 2007
$ws1.Range("A4").Value = 2008
$ws1.Range("A5").Value = 2009
$ws1.Range("A6").Value = 2010

$ws2 = $wb.Worksheets.Item("Page2")
$ws2.Range("A3").Value = 2007
$ws2.Range("A4").Value = 2008
$ws2.Range("A5").Value = 2009
$ws2.Range("A6").Value = 2010

# Make Page1 the active/selected sheet (previously Page2 was active)
$ws1.Activate()
$ws1.Range("A3:A6").Select()
